$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 12939.823
$ws.Range("I64").Value = 8998.799999999999
$ws.Range("K64").Value = 8998.799999999999
$ws.Range("M64").Value = -8750.799999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 12939.823
$ws.Range("I67").Value = 8998.799999999999
$ws.Range("K67").Value = 8998.799999999999
$ws.Range("M67").Value = -8140.799999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6655.8335
$ws.Range("I116").Value = 6384.1665
$ws.Range("J116").Value = 6791.6665
$ws.Range("K116").Value = 6384.1665
$ws.Range("L116").Value = 6791.6665
$ws.Range("M116").Value = -2942.1665
$ws.Range("N116").Value = -13675.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3936.6206
$ws.Range("I32").Value = 3283.5715
$ws.Range("K32").Value = 3283.5715
$ws.Range("M32").Value = -2996.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9966.333000000001
$ws.Range("I61").Value = 11449.5
$ws.Range("K61").Value = 11449.5
$ws.Range("M61").Value = -11237.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5650.2666
$ws.Range("I63").Value = 1499.75
$ws.Range("J63").Value = 7159.5454
$ws.Range("K63").Value = 1499.75
$ws.Range("L63").Value = 7159.5454
$ws.Range("M63").Value = -813.75
$ws.Range("N63").Value = -8531.545399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 5650.2666
$ws.Range("I66").Value = 1499.75
$ws.Range("J66").Value = 7159.5454
$ws.Range("K66").Value = 7498.75
$ws.Range("L66").Value = 35797.727
$ws.Range("M66").Value = -4066.75
$ws.Range("N66").Value = -42661.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 9966.333000000001
$ws.Range("I136").Value = 11449.5
$ws.Range("K136").Value = 34348.5
$ws.Range("M136").Value = -31798.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 8624.875
$ws.Range("J14").Value = 1099
$ws.Range("L14").Value = 1099
$ws.Range("N14").Value = -1443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 3997
$ws.Range("J23").Value = 3997
$ws.Range("L23").Value = 3997
$ws.Range("N23").Value = -4563

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2995.2354
$ws.Range("I86").Value = 3241.7144
$ws.Range("K86").Value = 3241.7144
$ws.Range("M86").Value = -2118.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2995.2354
$ws.Range("I89").Value = 3241.7144
$ws.Range("K89").Value = 16208.572
$ws.Range("M89").Value = -10592.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2041.5454
$ws.Range("J99").Value = 2337
$ws.Range("L99").Value = 2337
$ws.Range("N99").Value = -5333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 26845.182
$ws.Range("J12").Value = 32777.445
$ws.Range("L12").Value = 32777.445
$ws.Range("N12").Value = -33117.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8531
$ws.Range("I58").Value = 8039
$ws.Range("J58").Value = 10007
$ws.Range("K58").Value = 8039
$ws.Range("L58").Value = 10007
$ws.Range("M58").Value = -7836
$ws.Range("N58").Value = -10413

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4997.5
$ws.Range("I62").Value = 4997.5
$ws.Range("K62").Value = 4997.5
$ws.Range("M62").Value = -4373.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4997.5
$ws.Range("I65").Value = 4997.5
$ws.Range("K65").Value = 24987.5
$ws.Range("M65").Value = -21867.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 7326.75
$ws.Range("I69").Value = 7326.75
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 7326.75
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -6577.75
$ws.Range("N69").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 7326.75
$ws.Range("I72").Value = 7326.75
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 21980.25
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -18236.25
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6797.579
$ws.Range("J99").Value = 9338.833000000001
$ws.Range("L99").Value = 9338.833000000001
$ws.Range("N99").Value = -12334.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 460.09525
$ws.Range("I107").Value = 203.26315
$ws.Range("K107").Value = 203.26315
$ws.Range("M107").Value = 1716.73685

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6797.579
$ws.Range("J126").Value = 9338.833000000001
$ws.Range("L126").Value = 28016.499
$ws.Range("N126").Value = -32956.499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7564.7144
$ws.Range("I132").Value = 4388.75
$ws.Range("K132").Value = 13166.25
$ws.Range("M132").Value = -10636.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3318.2727
$ws.Range("I134").Value = 2889
$ws.Range("K134").Value = 8667
$ws.Range("M134").Value = -6132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8531
$ws.Range("I136").Value = 8039
$ws.Range("J136").Value = 10007
$ws.Range("K136").Value = 24117
$ws.Range("L136").Value = 30021
$ws.Range("M136").Value = -21567
$ws.Range("N136").Value = -35121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 3000
$ws.Range("I10").Value = 3000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -8861
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 850
$ws.Range("I92").Value = 650
$ws.Range("K92").Value = 1950
$ws.Range("M92").Value = -702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 999.6667
$ws.Range("I97").Value = 999.5
$ws.Range("K97").Value = 2998.5
$ws.Range("M97").Value = -2502.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 4163
$ws.Range("J39").Value = 4163
$ws.Range("L39").Value = 4163
$ws.Range("N39").Value = -5227

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2196.4
$ws.Range("I122").Value = 1993.5
$ws.Range("J122").Value = 3008
$ws.Range("K122").Value = 5980.5
$ws.Range("L122").Value = 9024
$ws.Range("M122").Value = -3530.5
$ws.Range("N122").Value = -13924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 974.3333
$ws.Range("I16").Value = 1071.6
$ws.Range("K16").Value = 1071.6
$ws.Range("M16").Value = -901.5999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4555.4443
$ws.Range("I68").Value = 2625
$ws.Range("K68").Value = 2625
$ws.Range("M68").Value = -1876

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4555.4443
$ws.Range("I71").Value = 2625
$ws.Range("K71").Value = 13125
$ws.Range("M71").Value = -9381

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2964.8333
$ws.Range("I93").Value = 2998
$ws.Range("J93").Value = 2931.6667
$ws.Range("K93").Value = 2998
$ws.Range("L93").Value = 2931.6667
$ws.Range("M93").Value = -1750
$ws.Range("N93").Value = -5427.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 3000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -3280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
